# Insert a new weekly data row for "Ciboulette" (Vega Modelo de Temuco) right
# after the existing row 400. This shifts all subsequent rows (old 401..436)
# down by one (to 402..437) and adds one brand-new record at row 401.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 401 - pushes everything below down by one.
$ws.Rows.Item(401).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(401, 1).Value = 10
$ws.Cells.Item(401, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(401, 3).Value = "La Araucanía"
$ws.Cells.Item(401, 4).Value = 45223
$ws.Cells.Item(401, 5).Value = 9
$ws.Cells.Item(401, 6).Value = 100112039
$ws.Cells.Item(401, 7).Value = "Ciboulette"
$ws.Cells.Item(401, 8).Value = "Sin especificar"
$ws.Cells.Item(401, 9).Value = "Primera"
$ws.Cells.Item(401, 10).Value = 30
$ws.Cells.Item(401, 11).Value = 7000
$ws.Cells.Item(401, 12).Value = 7000
$ws.Cells.Item(401, 13).Value = 7000
$ws.Cells.Item(401, 14).Value = "`$/docena de atados"
$ws.Cells.Item(401, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(401, 16).Value = 2333
$ws.Cells.Item(401, 17).Value = 3
$ws.Cells.Item(401, 18).Value = "Hortaliza"
